$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.045.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.49%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.565.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.00%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.26%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'208.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.10%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +0.85%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.36%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'22.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.67%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.18%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  +1.91%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  +0.55%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.567.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.22%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  +0.81%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D15").Value = "'27.029.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.46%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'61.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.57%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  +1.14%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'216.11"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Value = "'7.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.15%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'1.00"
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = "'  +2.45%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'9.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.12%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  -0.12%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'154.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.05%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  +0.02%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +1.06%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +1.90%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'1.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.36%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +1.68%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  +3.98%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  +0.38%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +4.65%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.428.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.91%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +12.57%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +2.10%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +2.73%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +1.44%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  +2.18%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.814"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.80%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +1.68%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.35%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +0.51%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -0.19%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'64.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.48%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'  -0.25%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.700.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.00%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  -1.05%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.0₆0102"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.26%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0519"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.52%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +0.35%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.48%  "
$ws.Range("E51").Style = "Normal"
